$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Double Trouble" bracket sheet: fill in the "." placeholders in the O and F
# columns with the actual round winners, copying the cell formatting (font /
# fill) from the matching D/E (or P) column so the new entries look like the
# rest of the bracket.
# ---------------------------------------------------------------------------

function Set-BracketCell($targetAddr, $value, $formatSourceAddr) {
    $ws.Range($targetAddr).Value = $value
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# Column O winners (copy style from the matching P-column entry)
Set-BracketCell "O2"  "Gorilla"             "P2"
Set-BracketCell "O4"  "Wolverine"           "P4"
Set-BracketCell "O6"  "Seahorse"            "P7"
Set-BracketCell "O8"  "Indian Cobra"        "P8"
Set-BracketCell "O10" "Sarus Crane"         "P10"
Set-BracketCell "O12" "Largetooth Sawfish"  "P12"
Set-BracketCell "O14" "Red-Footed Booby"    "P14"
Set-BracketCell "O16" "Speartooth Shark"    "P16"

# Column F winners (copy style from the matching D/E-column entry)
Set-BracketCell "F3"  "Pygmy Hog"     "D2"
Set-BracketCell "F7"  "Tree Hyrax"    "D8"
Set-BracketCell "F11" "Raccoon Dog"   "D12"
Set-BracketCell "F15" "Pudu"          "D16"
Set-BracketCell "F20" "Sloth Bear"    "D19"
Set-BracketCell "F24" "African Civet" "D25"
Set-BracketCell "F28" "Brown Hyena"   "D29"
Set-BracketCell "F32" "Tayra"         "D31"

$ws.Application.CutCopyMode = $false

# Row 12 grew taller (wrapped text) in the edited version.
$ws.Rows("12").RowHeight = 38

# Scroll position / selection as left by the editor.
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F32").Select()
